$d = $word.ActiveDocument

# --- Helper accented characters (avoid encoding issues in the source file) ---
$aAcute  = [char]0x00E1   # á
$eAcute  = [char]0x00E9   # é
$iAcute  = [char]0x00ED   # í
$cCedil  = [char]0x00E7   # ç
$oAcute  = [char]0x00F3   # ó

# ==========================================================================
# 1) Paragraph "O funcionário deverá informar...": drop the " (1,n)" marker
#    that follows "telefone". The paragraph is made up of five runs; Word's
#    re-layout on a text edit would normally coalesce runs that share
#    identical formatting, so temporary bookmarks are dropped at each
#    existing run boundary first to pin those splits in place.
# ==========================================================================

$rBoundary1 = $d.Content
$rBoundary1.Find.Execute("ser" + $aAcute + " gerada uma matr" + $iAcute + "cula", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundary1 = $rBoundary1.End

$rBoundary2 = $d.Content
$rBoundary2.Find.Execute("matr" + $iAcute + "cula, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundary2 = $rBoundary2.End

$rBoundary3 = $d.Content
$rBoundary3.Find.Execute("matr" + $iAcute + "cula, senha", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundary3 = $rBoundary3.End

$rBoundary4 = $d.Content
$rBoundary4.Find.Execute("matr" + $iAcute + "cula, senha e c" + $oAcute + "digo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundary4 = $rBoundary4.End

$d.Bookmarks.Add("zzz_split1", $d.Range($boundary1, $boundary1))
$d.Bookmarks.Add("zzz_split2", $d.Range($boundary2, $boundary2))
$d.Bookmarks.Add("zzz_split3", $d.Range($boundary3, $boundary3))
$d.Bookmarks.Add("zzz_split4", $d.Range($boundary4, $boundary4))

$rTelefone = $d.Content
$rTelefone.Find.Execute("Nome, telefone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterTelefone = $rTelefone.End

$rMarker = $d.Content
$rMarker.Find.Execute("Nome, telefone (1,n)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterMarker = $rMarker.End

$d.Range($posAfterTelefone, $posAfterMarker).Delete()

$d.Bookmarks.Item("zzz_split1").Delete()
$d.Bookmarks.Item("zzz_split2").Delete()
$d.Bookmarks.Item("zzz_split3").Delete()
$d.Bookmarks.Item("zzz_split4").Delete()

# ==========================================================================
# 2) Paragraph "O fornecedor terá os seguintes atributos...": drop the
#    " (1, n)" marker that follows "telefone" at the very end of the
#    paragraph (keeping the final period). Pin the existing "O fornecedor" /
#    " terá..." run boundary first so it isn't coalesced by the edit.
# ==========================================================================

$rFornBoundary = $d.Content
$rFornBoundary.Find.Execute("O fornecedor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fornBoundary = $rFornBoundary.End
$d.Bookmarks.Add("zzz_forn_split", $d.Range($fornBoundary, $fornBoundary))

$rFornTelefone = $d.Content
$rFornTelefone.Find.Execute("c" + $oAcute + "digo_fornecedor, nome, endere" + $cCedil + "o (cep, logradouro, numero, bairro, cidade, uf), telefone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterFornTelefone = $rFornTelefone.End

$rFornMarker = $d.Content
$rFornMarker.Find.Execute("c" + $oAcute + "digo_fornecedor, nome, endere" + $cCedil + "o (cep, logradouro, numero, bairro, cidade, uf), telefone (1, n).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterFornPeriod = $rFornMarker.End
$posAfterFornMarker = $posAfterFornPeriod - 1

$d.Range($posAfterFornTelefone, $posAfterFornMarker).Delete()

$d.Bookmarks.Item("zzz_forn_split").Delete()

# ==========================================================================
# 3) Relocate the "_GoBack" bookmark: it was sitting right after "campos" in
#    the "pessoa jurídica" paragraph; move it to sit right after the
#    (now-trimmed) "telefone" in the fornecedor paragraph, ahead of the
#    trailing ".".
# ==========================================================================

$d.Bookmarks.Item("_GoBack").Delete()

$rFornFinal = $d.Content
$rFornFinal.Find.Execute("c" + $oAcute + "digo_fornecedor, nome, endere" + $cCedil + "o (cep, logradouro, numero, bairro, cidade, uf), telefone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $rFornFinal.End

$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))
